$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.712.22"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.597.82"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'211.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  +0.69%  "
$ws.Range("E9").Value = "  +1.16%  "
$ws.Range("D10").Value = "'19.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.48%  "
$ws.Range("D11").Value = "'0.0841"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").Value = "1.822.71"
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").Value = "1.598.59"
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").Value = "'65.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.46%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.0₃0765"
$ws.Range("E17").Value = "  +5.38%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "26.666.74"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").Value = "'209.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.42%  "
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("D21").Value = "'7.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.94%  "
$ws.Range("E22").Value = "  +1.35%  "
$ws.Range("D23").Value = "'2.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("D24").Value = "'8.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.41%  "
$ws.Range("D25").Value = "'143.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.62%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("E28").Value = "  +0.38%  "
$ws.Range("D29").Value = "'15.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.78%  "
$ws.Range("D30").Value = "'0.0518"
$ws.Range("D30").Style = "Normal"
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("D32").Value = "'3.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.03%  "
$ws.Range("D33").Value = "'2.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.89%  "
$ws.Range("D34").Value = "1.286.67"
$ws.Range("E34").Value = "  +0.64%  "
$ws.Range("D35").Value = "'0.619"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.75%  "
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("E37").Value = "  +0.47%  "
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("E39").Value = "  +16.89%  "
$ws.Range("E40").Value = "  -1.17%  "
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "'2.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.783"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.36%  "
$ws.Range("D44").Value = "'63.19"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").Value = "1.735.18"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("D46").Value = "'91.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.73%  "
$ws.Range("E47").Value = "  -2.13%  "
$ws.Range("E48").Value = "  +0.25%  "
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("D51").Value = "'7.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.24%  "
